$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the _GoBack bookmark that used to sit right after "MEETING 1"
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 2. Collapse split runs back into single runs (these were apparently split
#    by spell-check / re-typing and the edit just re-merges them).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(": Rivera, Deiters and Team Adventuras", $false, $false, $false, $false, $false, $true, 1, $false, ": Rivera, Deiters and Team Adventuras", 2) | Out-Null

$d.Content.Find.Execute(":  8/30/2017 1730-1823", $false, $false, $false, $false, $false, $true, 1, $false, ":  8/30/2017 1730-1823", 2) | Out-Null

$d.Content.Find.Execute(": Dr.Moore, Rivera, Deiters and both GGC Maps Teams", $false, $false, $false, $false, $false, $true, 1, $false, ": Dr.Moore, Rivera, Deiters and both GGC Maps Teams", 2) | Out-Null

$d.Content.Find.Execute("Finds legend toggle feature clunky and difficult. Should move if needed.", $false, $false, $false, $false, $false, $true, 1, $false, "Finds legend toggle feature clunky and difficult. Should move if needed.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Append the MEETING 3 minutes after the last paragraph of the document
#    (right before the final section break).
# ---------------------------------------------------------------------------
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$insertionPara = $d.Paragraphs.Last
$insertionPara.Range.ListFormat.RemoveNumbers()
$insertionPara.Style = $d.Styles("Normal")

$newBodyXml = @'
<w:p/>
<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>MEETING 3</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">9/21/2017 </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>1640 - 1800</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Attending:</w:t></w:r><w:r><w:t xml:space="preserve"> David Rivera and Mike Deiters, Team A</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Dr. Moore unable to attend.</w:t></w:r></w:p>
<w:p><w:r><w:t>Bryan: Stated his work on events and how it functions. (Look at Get Help and Created by Tab for reference about displaying the same style at the current app &#8211; Dieters)</w:t></w:r></w:p>
<w:p><w:r><w:t>Matt: Divs will not work in svg, use groups instead.</w:t></w:r></w:p>
<w:p><w:r><w:t>Dieters: Asked how event function will be dynamic. (Using a server &#8211; Bryan)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
<w:p><w:r><w:t xml:space="preserve">Matt: States setting up tester units for app, and an onHover tooltip that displays name of current target, and if it can be added to function that changes svg colors on hover. (Look at Nav.scss file, poss nested in another class. Pull out of nested class if necessary. &#8211; Mike&amp;David) </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Popup.class inside </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">PopupGroup inside </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>content</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>. &#8211; Uses the css which is compiled from the layout scss.</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">Matt: Are there existing testing scripts for the app by prior groups? (2 or 3 but unsure their function/use due to Maggie </w:t></w:r><w:r><w:t>overseeing</w:t></w:r><w:r><w:t xml:space="preserve"> testing. Usually only made a test unit in response to a known bug to help fix it. &#8211; David)</w:t></w:r></w:p>
<w:p><w:r><w:t>Robert: Working on search suggestion, and made a Allrooms.txt to store all rooms on campus in an organized fashion to be used for search function. Asked about using datalist for search. (You can use datalist from HTML5. If using a dropdown, stylizing the select tag &#8216;sucks&#8217;. Check for custom dropdown style)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Google may use div with an unordered list.</w:t></w:r></w:p>
<w:p><w:r><w:t>Robert: Could not get app to recognize new functions in client.js, but it works fine in script.min.js where we shouldn&#8217;t be coding. What am I misunderstanding about the code? (Grunt is needed for app to read the client.js file, so please don&#8217;t code in script.min.js, and have Grunt running when testing your code.)</w:t></w:r></w:p>
<w:p><w:r><w:t>Look up Grunt plugin for atom to work with atom-live-server.</w:t></w:r></w:p>
<w:p><w:r><w:t>Do not list personal living quarters or faculty office names.</w:t></w:r></w:p>
'@

$xmlPackage = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$newBodyXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertionPara.Range.InsertXML($xmlPackage)

Write-Output "done"
